$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.723.20"
$ws.Range("E2").Value = "  -2.43%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.437.08"
$ws.Range("E3").Value = "  -5.14%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.96"
$ws.Range("E5").Value = "  -5.17%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.99"
$ws.Range("E6").Value = "  -3.99%  "

# Row 7
$ws.Range("E7").Value = "  -2.81%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.425.80"
$ws.Range("E8").Value = "  -5.12%  "

# Row 9
$ws.Range("E9").Value = "  +0.08%  "

# Row 10
$ws.Range("E10").Value = "  -6.74%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.611"
$ws.Range("E11").Value = "  -5.14%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "50.73"
$ws.Range("E12").Value = "  -4.84%  "

# Row 13
$ws.Range("E13").Value = "  -7.84%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.99"
$ws.Range("E14").Value = "  -5.61%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.974.96"
$ws.Range("E15").Value = "  -5.25%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "631.21"
$ws.Range("E16").Value = "  +3.98%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.539.90"
$ws.Range("E17").Value = "  -2.77%  "

# Row 18
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.120"
$ws.Range("E18").Value = "  -2.28%  "

# Row 19
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.427.16"
$ws.Range("E19").Value = "  -5.46%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.13"

# Row 21
$ws.Range("E21").Value = "  -5.52%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.932"
$ws.Range("E22").Value = "  -6.43%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.65"
$ws.Range("E23").Value = "  -2.70%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.24"
$ws.Range("E24").Value = "  -0.78%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.94"
$ws.Range("E25").Value = "  -5.45%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.22"
$ws.Range("E26").Value = "  -8.30%  "

# Row 27
$ws.Range("E27").Value = "  +1.96%  "

# Row 28
$ws.Range("E28").Value = "  -5.91%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.76"
$ws.Range("E29").Value = "  -7.81%  "

# Row 30
$ws.Range("E30").Value = "  -6.09%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.07"
$ws.Range("E31").Value = "  -4.98%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.11"
$ws.Range("E32").Value = "  -12.56%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.61"
$ws.Range("E33").Value = "  -7.92%  "

# Row 34
$ws.Range("E34").Value = "  -6.47%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "60.71"
$ws.Range("E35").Value = "  -4.09%  "

# Row 36
$ws.Range("E36").Value = "  -6.91%  "

# Row 37
$ws.Range("E37").Value = "  +0.02%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.634.49"
$ws.Range("E38").Value = "  -8.10%  "

# Row 39
$ws.Range("E39").Value = "  -12.75%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "482.99"
$ws.Range("E40").Value = "  -6.48%  "

# Row 41
$ws.Range("E41").Value = "  -7.22%  "

# Row 42
$ws.Range("E42").Value = "  -3.33%  "

# Row 43
$ws.Range("E43").Value = "  -6.55%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.132"
$ws.Range("E44").Value = "  -3.26%  "

# Row 45
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "33.84"
$ws.Range("E45").Value = "  -7.50%  "

# Row 46
$ws.Range("B46").Value = "CoreDAO"
$ws.Range("C46").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.46"
$ws.Range("E46").Value = "  +65.64%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0434"
$ws.Range("E47").Value = "  -5.89%  "

# Row 48
$ws.Range("E48").Value = "  -5.55%  "

# Row 49
$ws.Range("E49").Value = "  -4.75%  "

# Row 50
$ws.Range("E50").Value = "  -5.19%  "

# Row 51
$ws.Range("E51").Value = "  -0.30%  "
